$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("instructions")
Write-Host $ws.Name
